$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Title cell A1 gets the label "Год выпуска"
$ws.Range("A1").Value = "Год выпуска"

# 2. Column A gets a bit wider (14.4 -> ~15.6 chars)
$ws.Columns.Item(1).ColumnWidth = 14.8

# 3. The numeric comparison matrix (B2:G7) plus the priority-vector column (H2:H7)
#    become centered (horizontally AND vertically) in addition to the wrap already present.
$data = $ws.Range("B2:H7")
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108
$data.WrapText = $true

Write-Host "done"
